# Auto-generated edit script: updates market-price derived columns (H-N)
# across multiple sheets, matching the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 831.6667
$ws.Range("I9").Value = 247.5
$ws.Range("K9").Value = 247.5
$ws.Range("M9").Value = -78.5
# Row 64
$ws.Range("H64").Value = 200005420
$ws.Range("I64").Value = 333336400
$ws.Range("J64").Value = 8999.5
$ws.Range("K64").Value = 333336400
$ws.Range("L64").Value = 8999.5
$ws.Range("M64").Value = -333336152
$ws.Range("N64").Value = -9495.5
# Row 67
$ws.Range("H67").Value = 200005420
$ws.Range("I67").Value = 333336400
$ws.Range("J67").Value = 8999.5
$ws.Range("K67").Value = 333336400
$ws.Range("L67").Value = 8999.5
$ws.Range("M67").Value = -333335542
$ws.Range("N67").Value = -10715.5
# Row 69
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()
# Row 72
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()
# Row 100
$ws.Range("H100").Value = 2223.1304
$ws.Range("I100").Value = 1484.6666
$ws.Range("K100").Value = 1484.6666
$ws.Range("M100").Value = -943.6666
# Row 106
$ws.Range("H106").Value = 76924510
$ws.Range("I106").Value = 90910430
$ws.Range("K106").Value = 90910430
$ws.Range("M106").Value = -90909799
# Row 112
$ws.Range("H112").Value = 4775.4424
$ws.Range("I112").Value = 771
$ws.Range("J112").Value = 4935.62
$ws.Range("K112").Value = 2313
$ws.Range("L112").Value = 14806.86
$ws.Range("M112").Value = -1205
$ws.Range("N112").Value = -17022.86
# Row 129
$ws.Range("H129").Value = 1110.4
$ws.Range("I129").Value = 832.5
$ws.Range("J129").Value = 2222
$ws.Range("K129").Value = 2497.5
$ws.Range("L129").Value = 6666
$ws.Range("M129").Value = 2502.5
$ws.Range("N129").Value = -16666
# Row 138
$ws.Range("H138").Value = 4475.3335
$ws.Range("I138").Value = 1121.4286
$ws.Range("J138").Value = 7158.457
$ws.Range("K138").Value = 3364.2858
$ws.Range("L138").Value = 21475.371
$ws.Range("M138").Value = 1775.7142
$ws.Range("N138").Value = -31755.371

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3736.02
$ws.Range("I32").Value = 3758.4285
$ws.Range("K32").Value = 3758.4285
$ws.Range("M32").Value = -3471.4285
# Row 61
$ws.Range("H61").Value = 7737.147
$ws.Range("I61").Value = 4166.273
$ws.Range("K61").Value = 4166.273
$ws.Range("M61").Value = -3954.273
# Row 64
$ws.Range("H64").Value = 32332.572
$ws.Range("J64").Value = 32332.572
$ws.Range("L64").Value = 32332.572
$ws.Range("N64").Value = -32828.572
# Row 67
$ws.Range("H67").Value = 32332.572
$ws.Range("J67").Value = 32332.572
$ws.Range("L67").Value = 32332.572
$ws.Range("N67").Value = -34048.572
# Row 102
$ws.Range("H102").Value = 1019.8461
$ws.Range("I102").Value = 1040.3334
$ws.Range("J102").Value = 973.75
$ws.Range("K102").Value = 1040.3334
$ws.Range("L102").Value = 973.75
$ws.Range("M102").Value = 581.6666
$ws.Range("N102").Value = -4217.75
# Row 136
$ws.Range("H136").Value = 7737.147
$ws.Range("I136").Value = 4166.273
$ws.Range("K136").Value = 12498.819
$ws.Range("M136").Value = -9948.819

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 3257.35
$ws.Range("I94").Value = 779.0769
$ws.Range("J94").Value = 7859.857
$ws.Range("K94").Value = 779.0769
$ws.Range("L94").Value = 7859.857
$ws.Range("M94").Value = -328.0769
$ws.Range("N94").Value = -8761.857
# Row 134
$ws.Range("H134").Value = 5337.1133
$ws.Range("I134").Value = 1772
$ws.Range("K134").Value = 5316
$ws.Range("M134").Value = -2781

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5379.079
$ws.Range("I31").Value = 2281.6667
$ws.Range("J31").Value = 11697.8
$ws.Range("K31").Value = 2281.6667
$ws.Range("L31").Value = 11697.8
$ws.Range("M31").Value = -1986.6667
$ws.Range("N31").Value = -12287.8
# Row 34
$ws.Range("H34").Value = 5379.079
$ws.Range("I34").Value = 2281.6667
$ws.Range("J34").Value = 11697.8
$ws.Range("K34").Value = 2281.6667
$ws.Range("L34").Value = 11697.8
$ws.Range("M34").Value = -2079.6667
$ws.Range("N34").Value = -12101.8
# Row 127
$ws.Range("H127").Value = 50249
$ws.Range("J127").Value = 50249
$ws.Range("L127").Value = 50249
$ws.Range("N127").Value = -60169
# Row 132
$ws.Range("H132").Value = 3450.9844
$ws.Range("I132").Value = 1748.6111
$ws.Range("J132").Value = 12643.8
$ws.Range("K132").Value = 5245.8333
$ws.Range("L132").Value = 37931.39999999999
$ws.Range("M132").Value = -2715.8333
$ws.Range("N132").Value = -42991.39999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 35854732
$ws.Range("I4").Value = 49444444
$ws.Range("K4").Value = 148333332
$ws.Range("M4").Value = -148333220

$ws = $wb.Worksheets.Item("GSM")
# Row 47
$ws.Range("H47").Value = 33333.332
$ws.Range("J47").Value = 33333.332
$ws.Range("L47").Value = 33333.332
$ws.Range("N47").Value = -34469.332
# Row 97
$ws.Range("H97").Value = 1780.5
$ws.Range("J97").Value = 1707.6666
$ws.Range("L97").Value = 1707.6666
$ws.Range("N97").Value = -2699.6666
# Row 102
$ws.Range("H102").Value = 3010.9443
$ws.Range("I102").Value = 2830.3794
$ws.Range("J102").Value = 3759
$ws.Range("K102").Value = 2830.3794
$ws.Range("L102").Value = 3759
$ws.Range("M102").Value = -1208.3794
$ws.Range("N102").Value = -7003
# Row 122
$ws.Range("H122").Value = 2338328.2
$ws.Range("I122").Value = 2787006.8
$ws.Range("K122").Value = 8361020.399999999
$ws.Range("M122").Value = -8358570.399999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6852.706
$ws.Range("I7").Value = 5249.6
$ws.Range("K7").Value = 5249.6
$ws.Range("M7").Value = -5137.6
# Row 22
$ws.Range("H22").Value = 2654.7778
$ws.Range("I22").Value = 959.3333
$ws.Range("J22").Value = 3502.5
$ws.Range("K22").Value = 959.3333
$ws.Range("L22").Value = 3502.5
$ws.Range("M22").Value = -664.3333
$ws.Range("N22").Value = -4092.5
# Row 27
$ws.Range("H27").Value = 2654.7778
$ws.Range("I27").Value = 959.3333
$ws.Range("J27").Value = 3502.5
$ws.Range("K27").Value = 959.3333
$ws.Range("L27").Value = 3502.5
$ws.Range("M27").Value = -852.3333
$ws.Range("N27").Value = -3716.5
# Row 40
$ws.Range("H40").Value = 4195.613
$ws.Range("I40").Value = 2307.1428
$ws.Range("J40").Value = 8161.4
$ws.Range("K40").Value = 2307.1428
$ws.Range("L40").Value = 8161.4
$ws.Range("M40").Value = -2171.1428
# Row 126
$ws.Range("H126").Value = 6852.706
$ws.Range("I126").Value = 5249.6
$ws.Range("K126").Value = 15748.8
$ws.Range("M126").Value = -13278.8
# Row 138
$ws.Range("H138").Value = 75000
$ws.Range("J138").Value = 75000
$ws.Range("L138").Value = 75000
$ws.Range("N138").Value = -85280

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 25000
$ws.Range("J54").Value = 25000
$ws.Range("L54").Value = 25000
$ws.Range("N54").Value = -26040
# Row 126
$ws.Range("H126").Value = 3423.25
$ws.Range("I126").Value = 3167.1667
$ws.Range("K126").Value = 9501.500100000001
$ws.Range("M126").Value = -7031.500100000001
# Row 129
$ws.Range("H129").Value = 72952.336
$ws.Range("J129").Value = 72952.336
$ws.Range("L129").Value = 72952.336
$ws.Range("N129").Value = -82952.336
# Row 132
$ws.Range("H132").Value = 25020528
$ws.Range("I132").Value = 35723500
$ws.Range("J132").Value = 46925.832
$ws.Range("K132").Value = 107170500
$ws.Range("L132").Value = 140777.496
$ws.Range("M132").Value = -107167970
$ws.Range("N132").Value = -145837.496
# Row 133
$ws.Range("H133").Value = 137428.75
$ws.Range("J133").Value = 137428.75
$ws.Range("L133").Value = 137428.75
$ws.Range("N133").Value = -147548.75
